# Weekly crime-stats refresh: new data collected for the week of 3/18/2024-3/24/2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: bump issue number & reporting week ---
$ws.Range("A8").Value = "Volume 31   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/18/2024  Through  3/24/2024"

# --- Template cells already holding the correct "N/A" text + style we need to clone ---
# C14 = "0" (style 14), E14 = "***.*" (style 14)

# --- Cells that must become the text placeholder "0" (no data this period) ---
$ws.Range("C14").Copy() | Out-Null
$zeroTargets = @("D15","D27","C28","C29","D29","C30","D30","C33")
foreach ($cell in $zeroTargets) {
    $ws.Range($cell).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# --- Cells that must become the text placeholder "***.*" (pct change undefined) ---
$ws.Range("E14").Copy() | Out-Null
$naTargets = @("E15","E27","E29","E30")
foreach ($cell in $naTargets) {
    $ws.Range($cell).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# --- Updated numeric figures ---
$ws.Range("C15").Value = 2
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 2
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("H15").Value = 100
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I15").Value = 2
$ws.Range("I15").NumberFormat = '#,##0'
$ws.Range("K15").Value = -33.333333333333
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L15").Value = 0
$ws.Range("L15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C16").Value = 2
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 3
$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("E16").Value = -33.333333333333
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F16").Value = 5
$ws.Range("F16").NumberFormat = '#,##0'
$ws.Range("G16").Value = 6
$ws.Range("G16").NumberFormat = '#,##0'
$ws.Range("H16").Value = -16.666666666666
$ws.Range("H16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I16").Value = 15
$ws.Range("I16").NumberFormat = '#,##0'
$ws.Range("J16").Value = 18
$ws.Range("J16").NumberFormat = '#,##0'
$ws.Range("K16").Value = -16.666666666666
$ws.Range("K16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L16").Value = -11.764705882352
$ws.Range("L16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C17").Value = 5
$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 4
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("E17").Value = 25
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G17").Value = 16
$ws.Range("G17").NumberFormat = '#,##0'
$ws.Range("H17").Value = 18.75
$ws.Range("H17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I17").Value = 47
$ws.Range("I17").NumberFormat = '#,##0'
$ws.Range("J17").Value = 53
$ws.Range("J17").NumberFormat = '#,##0'
$ws.Range("K17").Value = -11.320754716981
$ws.Range("K17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L17").Value = -2.083333333333
$ws.Range("L17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D18").Value = 2
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("F18").Value = 2
$ws.Range("F18").NumberFormat = '#,##0'
$ws.Range("G18").Value = 3
$ws.Range("G18").NumberFormat = '#,##0'
$ws.Range("H18").Value = -33.333333333333
$ws.Range("H18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J18").Value = 9
$ws.Range("J18").NumberFormat = '#,##0'
$ws.Range("K18").Value = 33.333333333333
$ws.Range("K18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L18").Value = -42.857142857142
$ws.Range("L18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C19").Value = 6
$ws.Range("C19").NumberFormat = '#,##0'
$ws.Range("D19").Value = 9
$ws.Range("D19").NumberFormat = '#,##0'
$ws.Range("E19").Value = -33.333333333333
$ws.Range("E19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F19").Value = 25
$ws.Range("F19").NumberFormat = '#,##0'
$ws.Range("G19").Value = 27
$ws.Range("G19").NumberFormat = '#,##0'
$ws.Range("H19").Value = -7.407407407407
$ws.Range("H19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I19").Value = 81
$ws.Range("I19").NumberFormat = '#,##0'
$ws.Range("J19").Value = 69
$ws.Range("J19").NumberFormat = '#,##0'
$ws.Range("K19").Value = 17.391304347826
$ws.Range("K19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = -26.363636363636
$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F20").Value = 3
$ws.Range("F20").NumberFormat = '#,##0'
$ws.Range("G20").Value = 7
$ws.Range("G20").NumberFormat = '#,##0'
$ws.Range("H20").Value = -57.142857142857
$ws.Range("H20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I20").Value = 6
$ws.Range("I20").NumberFormat = '#,##0'
$ws.Range("J20").Value = 19
$ws.Range("J20").NumberFormat = '#,##0'
$ws.Range("K20").Value = -68.421052631578
$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L20").Value = -68.421052631578
$ws.Range("L20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D21").Value = 21
$ws.Range("D21").NumberFormat = '#,##0'
$ws.Range("E21").Value = -23.809523809523
$ws.Range("E21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("F21").Value = 56
$ws.Range("F21").NumberFormat = '#,##0'
$ws.Range("G21").Value = 60
$ws.Range("G21").NumberFormat = '#,##0'
$ws.Range("H21").Value = -6.666666666666
$ws.Range("H21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("I21").Value = 163
$ws.Range("I21").NumberFormat = '#,##0'
$ws.Range("J21").Value = 171
$ws.Range("J21").NumberFormat = '#,##0'
$ws.Range("K21").Value = -4.678362573099
$ws.Range("K21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("L21").Value = -24.884792626728
$ws.Range("L21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("D23").Value = 3
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G23").Value = 3
$ws.Range("G23").NumberFormat = '#,##0'
$ws.Range("H23").Value = -100
$ws.Range("H23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J23").Value = 6
$ws.Range("J23").NumberFormat = '#,##0'
$ws.Range("K23").Value = -83.333333333333
$ws.Range("K23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L23").Value = -83.333333333333
$ws.Range("L23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C24").Value = 25
$ws.Range("C24").NumberFormat = '#,##0'
$ws.Range("D24").Value = 20
$ws.Range("D24").NumberFormat = '#,##0'
$ws.Range("E24").Value = 25
$ws.Range("E24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F24").Value = 128
$ws.Range("F24").NumberFormat = '#,##0'
$ws.Range("G24").Value = 76
$ws.Range("G24").NumberFormat = '#,##0'
$ws.Range("H24").Value = 68.421052631578
$ws.Range("H24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I24").Value = 374
$ws.Range("I24").NumberFormat = '#,##0'
$ws.Range("J24").Value = 260
$ws.Range("J24").NumberFormat = '#,##0'
$ws.Range("K24").Value = 43.846153846153
$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L24").Value = 40.601503759398
$ws.Range("L24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C25").Value = 11
$ws.Range("C25").NumberFormat = '#,##0'
$ws.Range("D25").Value = 14
$ws.Range("D25").NumberFormat = '#,##0'
$ws.Range("E25").Value = -21.428571428571
$ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F25").Value = 93
$ws.Range("F25").NumberFormat = '#,##0'
$ws.Range("G25").Value = 46
$ws.Range("G25").NumberFormat = '#,##0'
$ws.Range("H25").Value = 102.173913043478
$ws.Range("H25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I25").Value = 258
$ws.Range("I25").NumberFormat = '#,##0'
$ws.Range("J25").Value = 164
$ws.Range("J25").NumberFormat = '#,##0'
$ws.Range("K25").Value = 57.317073170731
$ws.Range("K25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L25").Value = 85.611510791366
$ws.Range("L25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C26").Value = 11
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 10
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = 10
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F26").Value = 56
$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("H26").Value = 51.351351351351
$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I26").Value = 140
$ws.Range("I26").NumberFormat = '#,##0'
$ws.Range("J26").Value = 115
$ws.Range("J26").NumberFormat = '#,##0'
$ws.Range("K26").Value = 21.739130434782
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = 48.936170212766
$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C27").Value = 2
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("F27").Value = 2
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("H27").Value = 100
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I27").Value = 4
$ws.Range("I27").NumberFormat = '#,##0'
$ws.Range("K27").Value = -33.333333333333
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L27").Value = 100
$ws.Range("L27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G28").Value = 9
$ws.Range("G28").NumberFormat = '#,##0'
$ws.Range("H28").Value = -66.666666666666
$ws.Range("H28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J28").Value = 17
$ws.Range("J28").NumberFormat = '#,##0'
$ws.Range("K28").Value = -47.058823529411
$ws.Range("K28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L28").Value = -18.181818181818
$ws.Range("L28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L29").Value = -66.666666666666
$ws.Range("L29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L30").Value = -66.666666666666
$ws.Range("L30").NumberFormat = '#,##0.0;"-"#,##0.0'
